$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Crs")

# The CRS tracker's Status column (D2:D18) is driven off a single shared
# string; the reviewer flipped the outstanding "Waiting Response" rows to
# "Approved" (D2:D18 all shared that string, so updating the column
# re-points every row at the new text).
$ws.Range("D2:D18").Value = "Approved"

# Leave the selection where the editor apparently finished up.
$ws.Range("E14").Select()
